$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 585.5263
$ws.Range("I28").Value = 366.69232
$ws.Range("K28").Value = 366.69232
$ws.Range("M28").Value = 118.30768
$ws.Range("H86").Value = 3988.7932
$ws.Range("I86").Value = 2555.7144
$ws.Range("K86").Value = 2555.7144
$ws.Range("M86").Value = -1432.7144
$ws.Range("H89").Value = 3988.7932
$ws.Range("I89").Value = 2555.7144
$ws.Range("K89").Value = 12778.572
$ws.Range("M89").Value = -7162.572
$ws.Range("H106").Value = 3998.5715
$ws.Range("I106").Value = 3998.5715
$ws.Range("K106").Value = 3998.5715
$ws.Range("M106").Value = -3367.5715
$ws.Range("H112").Value = 6744.913
$ws.Range("J112").Value = 7438.1953
$ws.Range("L112").Value = 22314.5859
$ws.Range("N112").Value = -24530.5859
$ws.Range("H127").Value = 1682.8462
$ws.Range("I127").Value = 1682.8462
$ws.Range("K127").Value = 5048.5386
$ws.Range("M127").Value = -88.53859999999986
$ws.Range("H129").Value = 1810.7333
$ws.Range("I129").Value = 1570.0769
$ws.Range("K129").Value = 4710.2307
$ws.Range("M129").Value = 289.7692999999999
$ws.Range("H137").Value = 7105.909
$ws.Range("J137").Value = 18164.834
$ws.Range("L137").Value = 54494.50199999999
$ws.Range("N137").Value = -59594.50199999999
$ws.Range("H138").Value = 2880.2188
$ws.Range("I138").Value = 1408.35
$ws.Range("J138").Value = 5333.3335
$ws.Range("K138").Value = 4225.049999999999
$ws.Range("L138").Value = 16000.0005
$ws.Range("M138").Value = 914.9500000000007
$ws.Range("N138").Value = -26280.0005
$ws.Range("H141").Value = 2848.1428
$ws.Range("I141").Value = 2844.5833
$ws.Range("J141").Value = 2869.5
$ws.Range("K141").Value = 8533.749899999999
$ws.Range("L141").Value = 8608.5
$ws.Range("M141").Value = -3353.749899999999
$ws.Range("N141").Value = -18968.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23189.42
$ws.Range("J32").Value = 4999
$ws.Range("L32").Value = 4999
$ws.Range("N32").Value = -5573
$ws.Range("H122").Value = 3459.9092
$ws.Range("J122").Value = 4146.4546
$ws.Range("L122").Value = 12439.3638
$ws.Range("N122").Value = -17339.3638
$ws.Range("H132").Value = 3340166.5
$ws.Range("I132").Value = 6671333
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 20013999
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -20011469
$ws.Range("N132").Value = -32060

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 850.6
$ws.Range("I64").Value = 288.25
$ws.Range("K64").Value = 288.25
$ws.Range("M64").Value = -63.25
$ws.Range("H67").Value = 850.6
$ws.Range("I67").Value = 288.25
$ws.Range("K67").Value = 288.25
$ws.Range("M67").Value = 491.75
$ws.Range("H96").Value = 25214
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492
$ws.Range("H134").Value = 1895170
$ws.Range("I134").Value = 1834238
$ws.Range("K134").Value = 5502714
$ws.Range("M134").Value = -5500179

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 415.42856
$ws.Range("I7").Value = 468.7143
$ws.Range("K7").Value = 468.7143
$ws.Range("M7").Value = -355.7143
$ws.Range("H31").Value = 20665.074
$ws.Range("I31").Value = 9307
$ws.Range("J31").Value = 28473.75
$ws.Range("K31").Value = 9307
$ws.Range("L31").Value = 28473.75
$ws.Range("M31").Value = -9012
$ws.Range("N31").Value = -29063.75
$ws.Range("H34").Value = 20665.074
$ws.Range("I34").Value = 9307
$ws.Range("J34").Value = 28473.75
$ws.Range("K34").Value = 9307
$ws.Range("L34").Value = 28473.75
$ws.Range("M34").Value = -9105
$ws.Range("N34").Value = -28877.75
$ws.Range("H45").Value = 444
$ws.Range("I45").Value = 444
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 444
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 149
$ws.Range("H52").Value = 99990
$ws.Range("J52").Value = 99990
$ws.Range("L52").Value = 99990
$ws.Range("N52").Value = -100578
$ws.Range("H58").Value = 1033157.44
$ws.Range("J58").Value = 4468.6665
$ws.Range("L58").Value = 4468.6665
$ws.Range("N58").Value = -4874.6665
$ws.Range("H62").Value = 9999.6
$ws.Range("I62").Value = 4998
$ws.Range("K62").Value = 4998
$ws.Range("M62").Value = -4374
$ws.Range("H65").Value = 9999.6
$ws.Range("I65").Value = 4998
$ws.Range("K65").Value = 24990
$ws.Range("M65").Value = -21870
$ws.Range("H74").Value = 65156.5
$ws.Range("J74").Value = 65156.5
$ws.Range("L74").Value = 65156.5
$ws.Range("N74").Value = -66904.5
$ws.Range("H77").Value = 65156.5
$ws.Range("J77").Value = 65156.5
$ws.Range("L77").Value = 195469.5
$ws.Range("N77").Value = -204205.5
$ws.Range("H93").Value = 36592
$ws.Range("I93").Value = 39888
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 39888
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -38016
$ws.Range("N93").Value = -33744
$ws.Range("H99").Value = 3971.25
$ws.Range("I99").Value = 3895
$ws.Range("K99").Value = 3895
$ws.Range("M99").Value = -2397
$ws.Range("H126").Value = 3971.25
$ws.Range("I126").Value = 3895
$ws.Range("K126").Value = 11685
$ws.Range("M126").Value = -9215
$ws.Range("H132").Value = 3610.5557
$ws.Range("I132").Value = 3030.5
$ws.Range("J132").Value = 4074.6
$ws.Range("K132").Value = 9091.5
$ws.Range("L132").Value = 12223.8
$ws.Range("M132").Value = -6561.5
$ws.Range("N132").Value = -17283.8
$ws.Range("H136").Value = 1033157.44
$ws.Range("J136").Value = 4468.6665
$ws.Range("L136").Value = 13405.9995
$ws.Range("N136").Value = -18505.9995
$ws.Range("N45").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8600.454
$ws.Range("I56").Value = 8600.454
$ws.Range("K56").Value = 8600.454
$ws.Range("M56").Value = -8070.454
$ws.Range("H132").Value = 668.3333
$ws.Range("J132").Value = 865
$ws.Range("L132").Value = 7785
$ws.Range("N132").Value = -12845

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4151.5454
$ws.Range("I70").Value = 3897.4285
$ws.Range("K70").Value = 3897.4285
$ws.Range("M70").Value = -3627.4285
$ws.Range("H73").Value = 4151.5454
$ws.Range("I73").Value = 3897.4285
$ws.Range("K73").Value = 3897.4285
$ws.Range("M73").Value = -2961.4285
$ws.Range("H102").Value = 2997.8286
$ws.Range("I102").Value = 1985.9131
$ws.Range("J102").Value = 4937.3335
$ws.Range("K102").Value = 1985.9131
$ws.Range("L102").Value = 4937.3335
$ws.Range("M102").Value = -363.9131
$ws.Range("N102").Value = -8181.3335
$ws.Range("H122").Value = 115167
$ws.Range("J122").Value = 5008
$ws.Range("L122").Value = 15024
$ws.Range("N122").Value = -19924
$ws.Range("H132").Value = 4810.12
$ws.Range("I132").Value = 4003.8
$ws.Range("J132").Value = 5347.6665
$ws.Range("K132").Value = 12011.4
$ws.Range("L132").Value = 16042.9995
$ws.Range("M132").Value = -9481.400000000001
$ws.Range("N132").Value = -21102.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 11666.667
$ws.Range("J23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10460
$ws.Range("H40").Value = 2497.0952
$ws.Range("I40").Value = 2496.95
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 2496.95
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -2360.95
$ws.Range("N40").Value = -2772
$ws.Range("H122").Value = 4457.7417
$ws.Range("I122").Value = 4081.6667
$ws.Range("J122").Value = 6996.25
$ws.Range("K122").Value = 12245.0001
$ws.Range("L122").Value = 20988.75
$ws.Range("M122").Value = -9795.000100000001
$ws.Range("N122").Value = -25888.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1588.5625
$ws.Range("I122").Value = 1271.3043
$ws.Range("J122").Value = 2399.3333
$ws.Range("K122").Value = 3813.9129
$ws.Range("L122").Value = 7197.999899999999
$ws.Range("M122").Value = -1363.9129
$ws.Range("N122").Value = -12097.9999
